$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts existing rows 13:25 down to 14:26)
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the weekly update data
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C13").Value = "Los Lagos"
$ws.Range("D13").Value2 = 44803
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 100112013
$ws.Range("G13").Value = "Alcachofa"
$ws.Range("H13").Value = "Madrigal"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14500
$ws.Range("N13").Value = '$/caja 40 unidades'
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 362
$ws.Range("Q13").Value = 40
$ws.Range("R13").Value = "Hortaliza"
